$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Motores de aplicação"
$ws.Range("D3").Value = "-"

$ws.Range("C6").Value = "Motores de aplicação"
$ws.Range("D6").Value = "Usinagem"
